$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 10000
$ws.Range("J7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").Value = $null
$ws.Range("H14").Value = 10000
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").Value = $null
$ws.Range("H44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("N44").Value = $null

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 8500
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").Value = $null
$ws.Range("H14").Value = 2000
$ws.Range("J14").Value = 2000
$ws.Range("L14").Value = 2000
$ws.Range("N14").Value = -2350
$ws.Range("H20").Value = 8500
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").Value = $null
$ws.Range("H45").Value = 71429650
$ws.Range("I45").Value = 111112080
$ws.Range("J45").Value = 1282.8
$ws.Range("K45").Value = 111112080
$ws.Range("L45").Value = 1282.8
$ws.Range("M45").Value = -111111703
$ws.Range("N45").Value = -2036.8
$ws.Range("H122").Value = 1338.4242
$ws.Range("I122").Value = 1156.4615
$ws.Range("J122").Value = 2014.2858
$ws.Range("K122").Value = 3469.3845
$ws.Range("L122").Value = 6042.857400000001
$ws.Range("M122").Value = -1019.3845
$ws.Range("N122").Value = -10942.8574

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3129.35
$ws.Range("I16").Value = 1804.6923
$ws.Range("J16").Value = 5589.4287
$ws.Range("K16").Value = 1804.6923
$ws.Range("L16").Value = 5589.4287
$ws.Range("M16").Value = -1517.6923
$ws.Range("N16").Value = -6163.4287
$ws.Range("H19").Value = 206.1
$ws.Range("I19").Value = 206.1
$ws.Range("K19").Value = 206.1
$ws.Range("M19").Value = -36.09999999999999
$ws.Range("H24").Value = 206.1
$ws.Range("I24").Value = 206.1
$ws.Range("K24").Value = 206.1
$ws.Range("M24").Value = -36.09999999999999
$ws.Range("H25").Value = 19166.666
$ws.Range("H113").Value = 3129.35
$ws.Range("I113").Value = 1804.6923
$ws.Range("J113").Value = 5589.4287
$ws.Range("K113").Value = 1804.6923
$ws.Range("L113").Value = 5589.4287
$ws.Range("M113").Value = 365.3077000000001
$ws.Range("N113").Value = -9929.4287
$ws.Range("H122").Value = 1516.5476
$ws.Range("I122").Value = 1100.75
$ws.Range("J122").Value = 1772.4231
$ws.Range("K122").Value = 3302.25
$ws.Range("L122").Value = 5317.2693
$ws.Range("M122").Value = -852.25
$ws.Range("N122").Value = -10217.2693

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").Value = $null
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").Value = $null
$ws.Range("H122").Value = 2084042.1
$ws.Range("I122").Value = 638.5
$ws.Range("J122").Value = 6250849.5
$ws.Range("K122").Value = 5746.5
$ws.Range("L122").Value = 56257645.5
$ws.Range("M122").Value = -3296.5
$ws.Range("N122").Value = -56262545.5
$ws.Range("H123").Value = 1670
$ws.Range("I123").Value = 1010
$ws.Range("K123").Value = 3030
$ws.Range("M123").Value = -580

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").Value = $null
$ws.Range("H52").Value = 15806.6
$ws.Range("J52").Value = 15806.6
$ws.Range("L52").Value = 15806.6
$ws.Range("N52").Value = -16324.6
$ws.Range("H122").Value = 2997.375
$ws.Range("I122").Value = 3395.8
$ws.Range("J122").Value = 2333.3333
$ws.Range("K122").Value = 10187.4
$ws.Range("L122").Value = 6999.999899999999
$ws.Range("M122").Value = -7737.400000000001
$ws.Range("N122").Value = -11899.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").Value = $null
$ws.Range("H30").Value = 18233.334
$ws.Range("I30").Value = 700
$ws.Range("J30").Value = 27000
$ws.Range("K30").Value = 700
$ws.Range("L30").Value = 27000
$ws.Range("M30").Value = -592
$ws.Range("N30").Value = -27216
$ws.Range("H122").Value = 3139
$ws.Range("I122").Value = 3029.5293
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 9088.5879
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -6638.5879
$ws.Range("N122").Value = -19900

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 70200
$ws.Range("J16").Value = 70200
$ws.Range("L16").Value = 70200
$ws.Range("N16").Value = -70784
$ws.Range("H100").Value = 11364204
$ws.Range("I100").Value = 30303196
$ws.Range("J100").Value = 808
$ws.Range("K100").Value = 60606392
$ws.Range("L100").Value = 1616
$ws.Range("M100").Value = -60605851
$ws.Range("N100").Value = -2698
$ws.Range("H107").Value = 276.33334
$ws.Range("I107").Value = 314.4
$ws.Range("J107").Value = 241.72728
$ws.Range("K107").Value = 943.1999999999999
$ws.Range("L107").Value = 725.18184
$ws.Range("M107").Value = 976.8000000000001
$ws.Range("N107").Value = -4565.18184
$ws.Range("H113").Value = 500.65518
$ws.Range("I113").Value = 240.57143
$ws.Range("J113").Value = 743.4
$ws.Range("K113").Value = 721.71429
$ws.Range("L113").Value = 2230.2
$ws.Range("M113").Value = 1448.28571
$ws.Range("N113").Value = -6570.2
$ws.Range("H122").Value = 1255.8235
$ws.Range("I122").Value = 1085.5714
$ws.Range("J122").Value = 1375
$ws.Range("K122").Value = 3256.7142
$ws.Range("L122").Value = 4125
$ws.Range("M122").Value = -806.7142000000003
$ws.Range("N122").Value = -9025
$ws.Range("H136").Value = 3031.9
$ws.Range("I136").Value = 954.6667
$ws.Range("J136").Value = 5470.391
$ws.Range("K136").Value = 2864.0001
$ws.Range("L136").Value = 16411.173
$ws.Range("M136").Value = -314.0001000000002
$ws.Range("N136").Value = -21511.173
